$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dry_Weather")

$ws.Range("D2").Value = '"dry_weather_hourly"'
$ws.Range("E2").Value = '"dry_weather_monthly"'
$ws.Range("G2").ClearContents()
